$d = $word.ActiveDocument

# Update the first list item's text.
$d.Content.Find.Execute(
    "Slow down speed to about half (or ¾ speed, experiment).",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Have ghost mode timer pause while in frightened mode.", 2)

# Remove the now-duplicate paragraph that originally held that text
# (the second occurrence — the first one is the paragraph we just edited,
# which also carries the _GoBack bookmark that must be preserved).
$matchCount = 0
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r`a") -eq "Have ghost mode timer pause while in frightened mode.") {
        $matchCount++
        if ($matchCount -eq 2) {
            $p.Range.Delete()
            break
        }
    }
}
